$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "#! END_ROW"
$ws.Range("F1").Value = "#! END_ROW true"

$ws.Range("F1").Select()
